$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.449.07"
$ws.Range("E2").Value = "  -2.69%  "

# Row 3
$ws.Range("D3").Value = "3.484.07"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.98"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.17"
$ws.Range("E6").Value = "  -4.73%  "
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").Value = "  +3.68%  "
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.633"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.71"
$ws.Range("E11").Value = "  -4.87%  "
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "  -1.63%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.24"
$ws.Range("E13").Value = "  -2.25%  "
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "4.037.70"
$ws.Range("E14").Value = "  +0.42%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.75"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "3.482.30"
$ws.Range("E16").Value = "  +0.48%  "

# Row 17
$ws.Range("E17").Value = "  +0.37%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.08"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "65.410.41"
$ws.Range("E19").Value = "  -3.13%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.989"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "417.26"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "  +3.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "86.18"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "  -2.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.93"
$ws.Range("E25").Value = "  +9.61%  "
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.78"
$ws.Range("E26").Value = "  -9.83%  "
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.85"
$ws.Range("E27").Value = "  -2.94%  "
$ws.Range("D27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "  -3.51%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.08"
$ws.Range("E29").Value = "  +4.86%  "
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.26"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.53"
$ws.Range("E31").Value = "  -5.03%  "
$ws.Range("D31").Style = "Normal"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "614.23"
$ws.Range("E32").Value = "  -9.47%  "
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "  +0.36%  "

# Row 34
$ws.Range("E34").Value = "  -0.59%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.42"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("D35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "  +9.08%  "

# Row 37
$ws.Range("E37").Value = "  -0.07%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.40"
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "3.386.14"
$ws.Range("E39").Value = "  +11.03%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0785"
$ws.Range("E40").Value = "  -5.49%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.380"
$ws.Range("E41").Value = "  -5.97%  "
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "  +0.19%  "

# Row 43
$ws.Range("E43").Value = "  -5.98%  "

# Row 44
$ws.Range("E44").Value = "  -5.35%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.54"
$ws.Range("E45").Value = "  -8.81%  "
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.25"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0414"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "  -1.26%  "

# Row 49
$ws.Range("E49").Value = "  +1.91%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.48"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("D50").Style = "Normal"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.54"
$ws.Range("E51").Value = "  -1.96%  "
$ws.Range("D51").Style = "Normal"
